$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header labels for the new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header style used by the rest of row 1 (e.g. AC1) by copying its
# format (bold, centered/top aligned, thin box border) onto the new headers.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the team record values for every data row (2 through 48)
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 73   # AD
    $ws.Cells.Item($r, 31).Value = 89   # AE
    $ws.Cells.Item($r, 32).Value = 1    # AF
}
